$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting existing rows 39:170 down to 40:171.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new data record.
$ws.Range("A39").Value = 5
$ws.Range("B39").Value = "Macroferia Regional de Talca"
$ws.Range("C39").Value = "Maule"
$ws.Range("D39").Value = 44459
$ws.Range("E39").Value = 7
$ws.Range("F39").Value = 100114014
$ws.Range("G39").Value = "Betarraga"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 3000
$ws.Range("K39").Value = 700
$ws.Range("L39").Value = 700
$ws.Range("M39").Value = 700
$ws.Range("N39").Value = "$/paquete 5 unidades"
$ws.Range("O39").Value = "Región del Maule"
$ws.Range("P39").Value = 140
$ws.Range("Q39").Value = 5
$ws.Range("R39").Value = "Hortaliza"
